# Scheduled market-data refresh: update cached Universalis price/profit
# figures on the per-job Leve profit tables (one worksheet per Disciple
# of the Hand job). Only the price/profit columns (H-N) change; leve
# metadata (A-G) is untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 365.48
$ws.Range("I15").Value = 365.48
$ws.Range("K15").Value = 1096.44
$ws.Range("M15").Value = -927.4400000000001
$ws.Range("H55").Value = 496.33334
$ws.Range("I55").Value = 631.5
$ws.Range("K55").Value = 631.5
$ws.Range("M55").Value = -417.5
$ws.Range("H107").Value = 100003330
$ws.Range("I107").Value = 100003330
$ws.Range("K107").Value = 100003330
$ws.Range("M107").Value = -100001410

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1837.1875
$ws.Range("I61").Value = 1826.3334
$ws.Range("K61").Value = 1826.3334
$ws.Range("M61").Value = -1614.3334
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = $null
$ws.Range("H135").Value = 106653
$ws.Range("J135").Value = 106653
$ws.Range("L135").Value = 106653
$ws.Range("N135").Value = -116793
$ws.Range("H136").Value = 1837.1875
$ws.Range("I136").Value = 1826.3334
$ws.Range("K136").Value = 5479.0002
$ws.Range("M136").Value = -2929.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2531.7693
$ws.Range("I107").Value = 1940.4
$ws.Range("K107").Value = 1940.4
$ws.Range("M107").Value = -20.40000000000009
$ws.Range("H134").Value = 1444.7333
$ws.Range("I134").Value = 1190.7858
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 3572.3574
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -1037.3574
$ws.Range("N134").Value = -20070
$ws.Range("H135").Value = 78888.5
$ws.Range("J135").Value = 78888.5
$ws.Range("L135").Value = 78888.5
$ws.Range("N135").Value = -89028.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 12117.923
$ws.Range("I99").Value = 8838
$ws.Range("J99").Value = 14929.286
$ws.Range("K99").Value = 8838
$ws.Range("L99").Value = 14929.286
$ws.Range("M99").Value = -7340
$ws.Range("N99").Value = -17925.286
$ws.Range("H106").Value = 33249.5
$ws.Range("J106").Value = 33249.5
$ws.Range("L106").Value = 33249.5
$ws.Range("N106").Value = -35773.5
$ws.Range("H126").Value = 12117.923
$ws.Range("I126").Value = 8838
$ws.Range("J126").Value = 14929.286
$ws.Range("K126").Value = 26514
$ws.Range("L126").Value = 44787.858
$ws.Range("M126").Value = -24044
$ws.Range("N126").Value = -49727.858
$ws.Range("H134").Value = 2231.7778
$ws.Range("J134").Value = 3874
$ws.Range("L134").Value = 11622
$ws.Range("N134").Value = -16692

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 83626.336
$ws.Range("J2").Value = 537.3333
$ws.Range("L2").Value = 3223.9998
$ws.Range("N2").Value = -3449.9998
$ws.Range("H97").Value = 6333.8
$ws.Range("I97").Value = 5726.3335
$ws.Range("K97").Value = 17179.0005
$ws.Range("M97").Value = -16683.0005
$ws.Range("H98").Value = 4042
$ws.Range("I98").Value = 4265.6665
$ws.Range("J98").Value = 3874.25
$ws.Range("K98").Value = 12796.9995
$ws.Range("L98").Value = 11622.75
$ws.Range("M98").Value = -11298.9995
$ws.Range("N98").Value = -14618.75
$ws.Range("H131").Value = 1490
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 1490
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 4470
$ws.Range("M131").Value = $null
$ws.Range("N131").Value = -14550
$ws.Range("H132").Value = 3379.3076
$ws.Range("I132").Value = 2448.2727
$ws.Range("K132").Value = 22034.4543
$ws.Range("M132").Value = -19504.4543

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5132.5
$ws.Range("J80").Value = 5699.25
$ws.Range("L80").Value = 5699.25
$ws.Range("N80").Value = -7695.25
$ws.Range("H83").Value = 5132.5
$ws.Range("J83").Value = 5699.25
$ws.Range("L83").Value = 28496.25
$ws.Range("N83").Value = -38480.25
$ws.Range("H126").Value = 4694.25
$ws.Range("I126").Value = 4651.2
$ws.Range("K126").Value = 13953.6
$ws.Range("M126").Value = -11483.6
$ws.Range("H132").Value = 1933
$ws.Range("I132").Value = 1709.1428
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 5127.428400000001
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -2597.428400000001
$ws.Range("N132").Value = -15560

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 846.9048
$ws.Range("I16").Value = 771.5294
$ws.Range("K16").Value = 771.5294
$ws.Range("M16").Value = -601.5294
$ws.Range("H22").Value = 5849.6665
$ws.Range("I22").Value = 7933
$ws.Range("K22").Value = 7933
$ws.Range("M22").Value = -7638
$ws.Range("H27").Value = 5849.6665
$ws.Range("I27").Value = 7933
$ws.Range("K27").Value = 7933
$ws.Range("M27").Value = -7826
$ws.Range("H40").Value = 1673.0625
$ws.Range("I40").Value = 1617.9333
$ws.Range("K40").Value = 1617.9333
$ws.Range("M40").Value = -1481.9333
$ws.Range("H55").Value = 1497.5
$ws.Range("J55").Value = 1998
$ws.Range("L55").Value = 1998
$ws.Range("N55").Value = -2344
$ws.Range("H94").Value = 52500
$ws.Range("J94").Value = 52500
$ws.Range("L94").Value = 52500
$ws.Range("N94").Value = -53852
$ws.Range("H100").Value = 850
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").Value = $null
$ws.Range("H122").Value = 3719.4285
$ws.Range("I122").Value = 3719.4285
$ws.Range("K122").Value = 11158.2855
$ws.Range("M122").Value = -8708.2855
$ws.Range("H136").Value = 4325.6333
$ws.Range("I136").Value = 3904.8333
$ws.Range("K136").Value = 11714.4999
$ws.Range("M136").Value = -9164.499899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 22500
$ws.Range("I42").Value = 22500
$ws.Range("K42").Value = 22500
$ws.Range("M42").Value = -22122
$ws.Range("H43").Value = 14500.25
$ws.Range("I43").Value = 14500.25
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 14500.25
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -14351.25
$ws.Range("N43").Value = $null
$ws.Range("H136").Value = 1249
$ws.Range("I136").Value = 698.85
$ws.Range("K136").Value = 2096.55
$ws.Range("M136").Value = 453.4499999999998
